$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.399799346923828
$ws.Range("B1").Value = 1.59556519985199
$ws.Range("C1").Value = 1.985342741012573
$ws.Range("D1").Value = 2.062791109085083
$ws.Range("E1").Value = 1.589460134506226
